$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks on B2/C2 (login credentials no longer linked)
$ws.Hyperlinks.Delete()

# Update header row 1
$ws.Range("B1").Value = "User"

# Update row 2 credential values
$ws.Range("B2").Value = "usersetup"
$ws.Range("C2").Value = "b1f0rcE"

# Clear the "Invalid Credential" test block (row 4-5), keep formatting on row 4
$ws.Range("A4").ClearContents()
$ws.Range("B4:D4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("D5").ClearContents()

# Clear the "Forgot Password" test block (row 7-8), keep formatting on row 7
$ws.Range("A7").ClearContents()
$ws.Range("B7:C7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

# Update selected cell in the sheet view
$ws.Range("F5").Select()
